# Auto-generated Excel COM-interop script to apply the Aegis_Profits.xlsx market-data refresh
# described in the commit 'chore: update Sheets via scheduled runner'.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 92
$ws.Range("H92").Value = 1616.6
$ws.Range("I92").Value = 1771
$ws.Range("J92").Value = 999
$ws.Range("K92").Value = 1771
$ws.Range("L92").Value = 999
$ws.Range("M92").Value = -523
$ws.Range("N92").Value = -3495

# Row 116
$ws.Range("H116").Value = 2499.0476
$ws.Range("I116").Value = 2500
$ws.Range("J116").Value = 2498.889
$ws.Range("K116").Value = 2500
$ws.Range("L116").Value = 2498.889
$ws.Range("M116").Value = 942
$ws.Range("N116").Value = -9382.888999999999

# Row 137
$ws.Range("H137").Value = 1715.6666
$ws.Range("I137").Value = 1676.5
$ws.Range("J137").Value = 1833.1666
$ws.Range("K137").Value = 5029.5
$ws.Range("L137").Value = 5499.4998
$ws.Range("M137").Value = -2479.5
$ws.Range("N137").Value = -10599.4998

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 180.25
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 180.25
$ws.Range("K5").Value = 0
$ws.Range("L5").ClearContents()
$ws.Range("M5").Value = 180.25
$ws.Range("N5").Value = -404.25

# Row 25
$ws.Range("H25").Value = 750
$ws.Range("I25").Value = 750
$ws.Range("K25").Value = 750
$ws.Range("M25").Value = -348

# Row 61
$ws.Range("H61").Value = 1093.32
$ws.Range("I61").Value = 1014.4783
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1014.4783
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -802.4783
$ws.Range("N61").Value = -2424

# Row 95
$ws.Range("H95").Value = 23386.666
$ws.Range("J95").Value = 23386.666
$ws.Range("L95").Value = 23386.666
$ws.Range("N95").Value = -28878.666

# Row 103
$ws.Range("H103").Value = 45980
$ws.Range("J103").Value = 45980
$ws.Range("L103").Value = 45980
$ws.Range("N103").Value = -48324

# Row 109
$ws.Range("H109").Value = 35000
$ws.Range("J109").Value = 35000
$ws.Range("L109").Value = 35000
$ws.Range("N109").Value = -37774

# Row 132
$ws.Range("H132").Value = 2281.5264
$ws.Range("I132").Value = 1865.6875
$ws.Range("J132").Value = 4499.3335
$ws.Range("K132").Value = 5597.0625
$ws.Range("L132").Value = 13498.0005
$ws.Range("M132").Value = -3067.0625
$ws.Range("N132").Value = -18558.0005

# Row 136
$ws.Range("H136").Value = 1093.32
$ws.Range("I136").Value = 1014.4783
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 3043.4349
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -493.4349000000002
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 180.25
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 180.25
$ws.Range("K4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("M4").Value = 180.25
$ws.Range("N4").Value = -410.25

# Row 105
$ws.Range("H105").Value = 183579.1
$ws.Range("I105").Value = 101935.9
$ws.Range("J105").Value = 1000011
$ws.Range("K105").Value = 101935.9
$ws.Range("L105").Value = 1000011
$ws.Range("M105").Value = -100188.9
$ws.Range("N105").Value = -1003505

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 110
$ws.Range("I2").Value = 110
$ws.Range("K2").Value = 110
$ws.Range("M2").Value = 3

$ws = $wb.Worksheets.Item("CUL")
# Row 29
$ws.Range("H29").Value = 407.14285
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 407.14285
$ws.Range("K29").Value = 0
$ws.Range("L29").ClearContents()
$ws.Range("M29").Value = 1221.42855
$ws.Range("N29").Value = -1775.42855

# Row 46
$ws.Range("H46").Value = 56177.555
$ws.Range("I46").Value = 299
$ws.Range("J46").Value = 63162.375
$ws.Range("K46").Value = 897
$ws.Range("L46").Value = 189487.125
$ws.Range("M46").Value = -806
$ws.Range("N46").Value = -189669.125

# Row 60
$ws.Range("H60").Value = 436.66666
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

# Row 69
$ws.Range("H69").Value = 2750
$ws.Range("J69").Value = 2750
$ws.Range("L69").Value = 8250
$ws.Range("N69").Value = -9872

# Row 72
$ws.Range("H72").Value = 2750
$ws.Range("J72").Value = 2750
$ws.Range("L72").Value = 24750
$ws.Range("N72").Value = -32862

# Row 132
$ws.Range("H132").Value = 3375.6
$ws.Range("J132").Value = 3225.75
$ws.Range("L132").Value = 29031.75
$ws.Range("N132").Value = -34091.75

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 844551.7
$ws.Range("J5").Value = 16077.5
$ws.Range("L5").Value = 16077.5
$ws.Range("N5").Value = -16301.5

# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0

# Row 93
$ws.Range("H93").Value = 32000
$ws.Range("J93").Value = 32000
$ws.Range("L93").Value = 32000
$ws.Range("N93").Value = -35744

# Row 94
$ws.Range("H94").Value = 10448
$ws.Range("J94").Value = 10448
$ws.Range("L94").Value = 10448
$ws.Range("N94").Value = -11800

# Row 97
$ws.Range("H97").Value = 83335990
$ws.Range("I97").Value = 83335990
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 83335990
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -83335494

# Row 98
$ws.Range("H98").Value = 31643
$ws.Range("J98").Value = 31643
$ws.Range("L98").Value = 31643
$ws.Range("N98").Value = -37633

# Row 99
$ws.Range("H99").Value = 8155.1665
$ws.Range("I99").Value = 3786.2
$ws.Range("J99").Value = 30000
$ws.Range("K99").Value = 3786.2
$ws.Range("L99").Value = 30000
$ws.Range("M99").Value = -1540.2
$ws.Range("N99").Value = -34492

# Row 113
$ws.Range("H113").Value = 2883.5
$ws.Range("I113").Value = 3325.25
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 3325.25
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -1155.25
$ws.Range("N113").Value = -6340

# Row 126
$ws.Range("H126").Value = 2264670.8
$ws.Range("I126").Value = 2734.875
$ws.Range("J126").Value = 3269975.5
$ws.Range("K126").Value = 8204.625
$ws.Range("L126").Value = 9809926.5
$ws.Range("M126").Value = -5734.625
$ws.Range("N126").Value = -9814866.5

# Row 132
$ws.Range("H132").Value = 2456.074
$ws.Range("I132").Value = 1764
$ws.Range("K132").Value = 5292
$ws.Range("M132").Value = -2762

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 51697.35
$ws.Range("I40").Value = 143842.86
$ws.Range("J40").Value = 2080.5386
$ws.Range("K40").Value = 143842.86
$ws.Range("L40").Value = 2080.5386
$ws.Range("M40").Value = -143706.86
$ws.Range("N40").Value = -2352.5386

# Row 55
$ws.Range("H55").Value = 621.3158
$ws.Range("I55").Value = 330
$ws.Range("J55").Value = 791.25
$ws.Range("K55").Value = 330
$ws.Range("L55").Value = 791.25
$ws.Range("M55").Value = -157
$ws.Range("N55").Value = -1137.25

# Row 60
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").ClearContents()
$ws.Range("N60").Value = 0

# Row 61
$ws.Range("H61").Value = 2528.2307
$ws.Range("I61").Value = 1485
$ws.Range("J61").Value = 3422.4285
$ws.Range("K61").Value = 1485
$ws.Range("L61").Value = 3422.4285
$ws.Range("M61").Value = -1283
$ws.Range("N61").Value = -3826.4285

# Row 93
$ws.Range("H93").Value = 3491.8
$ws.Range("I93").Value = 3491.8
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 3491.8
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -2243.8

# Row 100
$ws.Range("H100").Value = 1819
$ws.Range("I100").Value = 1773.75
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1773.75
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1232.75
$ws.Range("N100").Value = -3082

# Row 104
$ws.Range("H104").Value = 21333.334
$ws.Range("J104").Value = 21333.334
$ws.Range("L104").Value = 21333.334
$ws.Range("N104").Value = -28321.334

# Row 105
$ws.Range("H105").Value = 46240
$ws.Range("J105").Value = 46240
$ws.Range("L105").Value = 46240
$ws.Range("N105").Value = -53228

# Row 113
$ws.Range("H113").Value = 2528.2307
$ws.Range("I113").Value = 1485
$ws.Range("J113").Value = 3422.4285
$ws.Range("K113").Value = 1485
$ws.Range("L113").Value = 3422.4285
$ws.Range("M113").Value = 685
$ws.Range("N113").Value = -7762.4285

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 3752166.8
$ws.Range("I5").Value = 5005000
$ws.Range("K5").Value = 5005000
$ws.Range("M5").Value = -5004888

# Row 96
$ws.Range("H96").Value = 142858540
$ws.Range("I96").Value = 250001550
$ws.Range("J96").Value = 1199.3334
$ws.Range("K96").Value = 250001550
$ws.Range("L96").Value = 1199.3334
$ws.Range("M96").Value = -250000177
$ws.Range("N96").Value = -3945.3334

# Row 97
$ws.Range("H97").Value = 40572
$ws.Range("J97").Value = 40572
$ws.Range("L97").Value = 40572
$ws.Range("N97").Value = -42554

# Row 98
$ws.Range("H98").Value = 29875
$ws.Range("J98").Value = 29875
$ws.Range("L98").Value = 29875
$ws.Range("N98").Value = -35865

# Row 100
$ws.Range("H100").Value = 168420.33
$ws.Range("I100").Value = 336333.34
$ws.Range("J100").Value = 507.33334
$ws.Range("K100").Value = 672666.6800000001
$ws.Range("L100").Value = 1014.66668
$ws.Range("M100").Value = -672125.6800000001
$ws.Range("N100").Value = -2096.66668

# Row 101
$ws.Range("H101").Value = 10000
$ws.Range("J101").Value = 10000
$ws.Range("L101").Value = 10000
$ws.Range("N101").Value = -16490

# Row 102
$ws.Range("H102").Value = 42158.5
$ws.Range("J102").Value = 42158.5
$ws.Range("L102").Value = 42158.5
$ws.Range("N102").Value = -48648.5

# Row 103
$ws.Range("H103").Value = 17000
$ws.Range("J103").Value = 17000
$ws.Range("L103").Value = 17000
$ws.Range("N103").Value = -19344

# Row 104
$ws.Range("H104").Value = 30000
$ws.Range("J104").Value = 30000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -36988

# Row 105
$ws.Range("H105").Value = 27500
$ws.Range("J105").Value = 27500
$ws.Range("L105").Value = 27500
$ws.Range("N105").Value = -34488

# Row 106
$ws.Range("H106").Value = 32000
$ws.Range("J106").Value = 32000
$ws.Range("L106").Value = 32000
$ws.Range("N106").Value = -34524

